# Apply the "Update on 2018-03-06, 支出生活费300" change to WanHaoBillDetails

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Rename the sheet from "Summary" to "第一学年" (also updates the
#    _FilterDatabase defined name automatically since it tracks the sheet)
$ws.Name = "第一学年"

# 2. Before overwriting row 29, stamp its current (blank-template) formatting
#    onto the 11 new rows that will follow it (rows 30-40), so they inherit
#    the same blank-row look the sheet already used for row 29.
$ws.Range("B29:G29").Copy()
$ws.Range("B30:G40").PasteSpecial(-4122)  # xlPasteFormats

# 3. Turn row 29 into a real data entry dated 2018-03-06 for a 300 生活费
#    expense, matching the formatting of the preceding data row (28).
$ws.Range("B28:G28").Copy()
$ws.Range("B29:G29").PasteSpecial(-4122)  # xlPasteFormats

$ws.Cells.Item(29, 2).Value2 = 27                     # B29 serial number
$ws.Cells.Item(29, 3).Value = "支出"                   # C29 category
$ws.Cells.Item(29, 4).Value2 = 300                     # D29 amount
$ws.Cells.Item(29, 5).Value2 = 43165                   # E29 date (2018-03-06)
$ws.Cells.Item(29, 6).Value = "生活费"                 # F29 expense type
$ws.Cells.Item(29, 7).Value = "生活费(3/6-3/15)"       # G29 note

# 4. Fill in the serial numbers for the new blank rows 30-40 (28..38) and
#    give column B the same "numbered row" look used elsewhere (style from B28).
$ws.Range("B28").Copy()
$ws.Range("B30:B40").PasteSpecial(-4122)  # xlPasteFormats

for ($i = 0; $i -lt 11; $i++) {
    $r = 30 + $i
    $ws.Cells.Item($r, 2).Value2 = 28 + $i
}

# 5. Scroll/select to reflect where the user ended up after the edit.
$excel.ActiveWindow.ScrollRow = 10
$ws.Range("E30").Select()

# Recalculate formulas (J3, K3, J9 depend on the new row's values)
$excel.Calculate()
